# dev-61 - fixed upload tt for many employees
# Add a new employee row (row 7) to the timetable, copying the visual
# formatting used by the existing rows above it, and move the active
# selection to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new employee "Смешнов Александр Николаевич" (ЗДМ),
#     number 26856, position "Продавец-кассир" ---

# Columns A:C use the same style as the header cells of other rows (style
# donor A6). Copy formatting first, then set the values.
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)

$ws.Range("A7").Value = 26856
$ws.Range("B7").Value = "Смешнов Александр Николаевич"
$ws.Range("C7").Value = "Продавец-кассир"

# Columns D:AG carry the work-schedule values, alternating between the
# "10:00-20:00" shift (white fill, style donor G6) and "В" day-off marker
# (green fill, style donor D6), matching the exact column layout below.

$shiftRanges = @("D7:F7", "J7:L7", "P7:R7", "W7:X7", "AB7:AD7")
$offRanges   = @("G7:I7", "M7:O7", "S7:V7", "Y7:AA7", "AE7:AG7")

foreach ($r in $shiftRanges) {
    $ws.Range("G6").Copy()
    $ws.Range($r).PasteSpecial(-4122)
    $ws.Range($r).Value = "10:00-20:00"
}

foreach ($r in $offRanges) {
    $ws.Range("D6").Copy()
    $ws.Range($r).PasteSpecial(-4122)
    $ws.Range($r).Value = "В"
}

$excel.CutCopyMode = $false

# --- Move the active selection to D7, matching the saved view state ---
[void]$ws.Range("D7").Select()
